$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "move_folder" header column (G), between publish_time and monetization ---
$ws.Range("G1").Value = "move_folder"

# --- Insert a new, empty row 2 beneath the header so the sheet's used
#     range / dimension grows to A1:H2 (matching autoFilter ref="A1:H2").
#     Do this before the bold-styling step below so the insert only has to
#     carry down the single style that already exists on H1. ---
$ws.Rows.Item(2).Insert()
$ws.Range("H2").ClearFormats()

# --- Make every header cell bold, reusing the existing bold style from H1
#     (copy/paste-special formats so all header cells share one style id) ---
$ws.Range("H1").Copy() | Out-Null
$ws.Range("A1:G1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Resize columns B, C and G to fit the wider "directory"/"title"/
#     "move_folder" data. ColumnWidth is expressed in "characters"; Excel
#     pads this by 5/6 of a character when writing the stored column
#     width, so subtract that out to land exactly on the target stored
#     widths (21, 13, 34). ---
$ws.Columns.Item(2).ColumnWidth = 20.16666666666667
$ws.Columns.Item(3).ColumnWidth = 12.16666666666667
$ws.Columns.Item(7).ColumnWidth = 33.16666666666667
